# Update GridSearchCV results (mean_fit_time, std_fit_time, mean_score_time, std_score_time)
# for rows 2-9 (index 0-7) to reflect the re-run results included with the
# new Unidade 05 clustering notebook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.09214415550231933
$ws.Range("C2").Value = 0.01068712843943281
$ws.Range("D2").Value = 0.5620160579681397
$ws.Range("E2").Value = 0.1879694053081425
$ws.Range("B3").Value = 0.08284244537353516
$ws.Range("C3").Value = 0.003576248863731651
$ws.Range("D3").Value = 0.4092443466186523
$ws.Range("E3").Value = 0.02078270958990371
$ws.Range("B4").Value = 0.08499932289123535
$ws.Range("C4").Value = 0.002901952515569321
$ws.Range("D4").Value = 0.3371630191802978
$ws.Range("E4").Value = 0.03410400802716913
$ws.Range("B5").Value = 0.0894345760345459
$ws.Range("C5").Value = 0.004878562157682315
$ws.Range("D5").Value = 0.3135540962219238
$ws.Range("E5").Value = 0.02533646317474679
$ws.Range("B6").Value = 0.08511719703674317
$ws.Range("C6").Value = 0.004477371024013886
$ws.Range("D6").Value = 0.4835751533508301
$ws.Range("E6").Value = 0.02100817168785565
$ws.Range("B7").Value = 0.1023220539093018
$ws.Range("C7").Value = 0.01750774223661531
$ws.Range("D7").Value = 0.5892225742340088
$ws.Range("E7").Value = 0.1146882500822285
$ws.Range("B8").Value = 0.08663473129272461
$ws.Range("C8").Value = 0.0050290070404196
$ws.Range("D8").Value = 0.3810472011566162
$ws.Range("E8").Value = 0.02514631095327748
$ws.Range("B9").Value = 0.08413500785827636
$ws.Range("C9").Value = 0.004418511037569464
$ws.Range("D9").Value = 0.344151496887207
$ws.Range("E9").Value = 0.01693298457819821
